$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# B8 was a numeric 2, now becomes the text "2+0.5"
$ws.Range("B8").Value = "2+0.5"

# C8 text gets appended with additional detail about VoiceImagesManager
$ws.Range("C8").Value = "Load data from Parse.com and display them on the ArticlesTableViewController. Build VoiceImagesManager to handle image"
